$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold + bordered, style index 1) by copying the
# formatting from H1 before writing the new text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J data columns for rows 2-35.
$data = @(
    @(7, 8),
    @(6, 8),
    @(1, 3),
    @(4, 6),
    @(4, 7),
    @(7, 7),
    @(4, 6),
    @(8, 8),
    @(8, 9),
    @(6, 6),
    @(4, 7),
    @(3, 5),
    @(5, 9),
    @(4, 6),
    @(5, 7),
    @(11, 11),
    @(4, 6),
    @(9, 9),
    @(1, 3),
    @(1, 4),
    @(4, 6),
    @(1, 3),
    @(6, 7),
    @(6, 8),
    @(6, 7),
    @(5, 6),
    @(7, 8),
    @(5, 6),
    @(1, 5),
    @(1, 4),
    @(4, 6),
    @(6, 7),
    @(2, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
